# Fruta / hortaliza, semanal
# A new week's price observation is inserted at the top of the data
# (row 11, right after the 9 "latest" rows above it), pushing the
# previously-existing rows 11-20 down to rows 12-21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 11, shifting rows 11:20 down to 12:21
# (mirrors Excel's Rows("11:11").Insert -> EntireRow.Insert with
# xlShiftDown, which also carries the D-column date style down).
$ws.Rows.Item(11).Insert()

# Populate the new row 11 with this week's observation.
$ws.Range("A11").Value = 7
$ws.Range("B11").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C11").Value = "Ñuble"
$ws.Range("D11").Value = 45096
$ws.Range("E11").Value = 16
$ws.Range("F11").Value = "Fruta"
$ws.Range("G11").Value = 100107
$ws.Range("H11").Value = "Otros"
$ws.Range("I11").Value = 100107011
$ws.Range("J11").Value = "Tuna"
$ws.Range("K11").Value = "Sin especificar"
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 50
$ws.Range("N11").Value = 23000
$ws.Range("O11").Value = 23000
$ws.Range("P11").Value = 23000
$ws.Range("Q11").Value = "$/caja 18 kilos"
$ws.Range("R11").Value = "Región Metropolitana"
$ws.Range("S11").Value = 1278
$ws.Range("T11").Value = 18
